{"js": "// Replace the label \"Talk\" with \"Direct instruction\" everywhere it appears\n// in the document body (it occurs as the whole contents of a table cell,\n// once per repeated sub-table \u2014 4 times total in this document).\nconst results = context.document.body.search(\"Talk\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // insertText(..., replace) overwrites only the matched range's text,\n  // leaving the run's existing formatting (font, size, color, ...) intact.\n  results.items[i].insertText(\"Direct instruction\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the label \"Talk\" with \"Direct instruction\" everywhere it appears\n# in the document's table (it occurs as the whole contents of a table cell,\n# once per repeated sub-table -- 4 times total in this document).\n$d = $word.ActiveDocument\n\nforeach ($table in $d.Tables) {\n    foreach ($row in $table.Rows) {\n        foreach ($cell in $row.Cells) {\n            # Cell text includes the trailing paragraph mark (Chr 13) and the\n            # end-of-cell marker (Chr 7); strip those before comparing.\n            $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n            if ($cellText -eq \"Talk\") {\n                # Assigning to .Text replaces only the text, leaving the\n                # existing run formatting (font, size, color, ...) intact.\n                $cell.Range.Text = \"Direct instruction\"\n            }\n        }\n    }\n}\n"}
